# InventoryCount.xlsx - "Getting ready for Lab 7 submission"
#
# 1. Insert a new "Item Selection" column between the existing
#    "ListBox Selection" (C) and "Action" (old D) columns.
# 2. Fill in the new column header + per-row item name.
# 3. Rename "Cookies" -> "Cookie" in the ListBox Selection column.
# 4. Tweak a couple of wording/height details on existing rows.
# 5. Add the new Lab 7 test rows (Click Add / Click Delete / tab order …)
#    at the bottom of the test-plan table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert the new ItemSelection column before column D -------------
$ws.Columns("D").Insert()
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# --- 2. New column header -------------------------------------------------
$ws.Range("D8").Value = "ItemSelection"

# --- 3. Rename Cookies -> Cookie, and populate the new ItemSelection column
$ws.Range("C10").Value = "Cookie"
$ws.Range("D10").Value = "Cookie"

$ws.Range("C11").Value = "Cookie"
$ws.Range("D11").Value = "Cookie"

$ws.Range("D12").Value = "Muffin"
$ws.Range("F12").Value = "The muffin's quantity will appear in the quantity textbox"
$ws.Rows("12").RowHeight = 60

# --- 4. New test rows for Lab 7 -------------------------------------------
$ws.Range("A20").Value = "Joe's Shop"
$ws.Range("E20").Value = "Click Add"
$ws.Range("F20").Value = "The item will be unselected, the quantity will be blank. No item will be selected in the listbox"
$ws.Range("G20").Value = "OK"
$ws.Rows("20").RowHeight = 90

$ws.Range("D23").Value = "Donut"
$ws.Range("D24").Value = "Donut"

$ws.Range("A21").Value = "Joe's Shop"
$ws.Range("E21").Value = "Click Save"
$ws.Range("F21").Value = "Error: Stock quantity between 0 and 300 must be entered, You must select a baked good."
$ws.Range("G21").Value = "OK"
$ws.Rows("21").RowHeight = 90

$ws.Range("A22").Value = "Joe's Shop"
$ws.Range("B22").Value = 1
$ws.Range("E22").Value = "Click Save"
$ws.Range("F22").Value = "Error: You must select a baked good."
$ws.Range("G22").Value = "OK"
$ws.Rows("22").RowHeight = 45

$ws.Range("A23").Value = "Joe's Shop"
$ws.Range("B23").Value = 1
$ws.Range("E23").Value = "Click Save"
$ws.Range("F23").Value = "Donut with quantity 1 added to listbox"
$ws.Range("G23").Value = "OK"
$ws.Rows("23").RowHeight = 45

$ws.Range("A24").Value = "Joe's Shop"
$ws.Range("B24").Value = 1
$ws.Range("E24").Value = "Click Delete"
$ws.Range("F24").Value = "Donut is deleted from the listbox"
$ws.Range("G24").Value = "OK"
$ws.Rows("24").RowHeight = 30

$ws.Range("A25").Value = "Joe's Shop"
$ws.Range("E25").Value = "Delete all items from the listbox and then click delete again"
$ws.Range("F25").Value = "Error: You must select an item to delete first"
$ws.Range("G25").Value = "OK"
$ws.Rows("25").RowHeight = 60

# --- 5. New row-7 instruction line (same row that already carried the ---
#        thick bottom border, just below the screen title) --------------
$ws.Range("A7").Value = "Tab order: Customer selection, listbox, item selection, quantity, Save, Add, Delete"
$ws.Range("A7").VerticalAlignment = -4108

# --- 6. Leave the selection on A7, scrolled back to the top of the sheet -
[void]$ws.Range("A7").Select()
